# This workbook ("Mandragora_Profits") is refreshed on a schedule by pulling
# current Universalis market-board prices per item (columns H:J) and
# recomputing the derived Leve buy-in / profit columns (K:N) for each class
# tab. This run updates the rows whose source prices moved since the last sync.

$wb = $excel.ActiveWorkbook

# ----- ALC sheet -----
$ws = $wb.Worksheets.Item("ALC")

# Row 55
$ws.Range("H55").Value = 390.46667
$ws.Range("I55").Value = 405.27274
$ws.Range("J55").Value = 349.75
$ws.Range("K55").Value = 405.27274
$ws.Range("L55").Value = 349.75
$ws.Range("M55").Value = -191.27274
$ws.Range("N55").Value = -777.75

# Row 58
$ws.Range("H58").Value = 2052.6667
$ws.Range("I58").Value = 285
$ws.Range("J58").Value = 2557.7144
$ws.Range("K58").Value = 855
$ws.Range("L58").Value = 7673.1432
$ws.Range("M58").Value = -705
$ws.Range("N58").Value = -7973.1432

# Row 69
$ws.Range("H69").Value = 6127.143
$ws.Range("I69").Value = 6986.6665
$ws.Range("J69").Value = 5482.5
$ws.Range("K69").Value = 20959.9995
$ws.Range("L69").Value = 16447.5
$ws.Range("M69").Value = -20085.9995
$ws.Range("N69").Value = -18195.5

# Row 72
$ws.Range("H72").Value = 6127.143
$ws.Range("I72").Value = 6986.6665
$ws.Range("J72").Value = 5482.5
$ws.Range("K72").Value = 62879.9985
$ws.Range("L72").Value = 49342.5
$ws.Range("M72").Value = -58511.9985
$ws.Range("N72").Value = -58078.5

# Row 103
$ws.Range("H103").Value = 1439.8
$ws.Range("I103").Value = 1033.3334
$ws.Range("J103").Value = 1710.7778
$ws.Range("K103").Value = 3100.0002
$ws.Range("L103").Value = 5132.3334
$ws.Range("M103").Value = -2514.0002
$ws.Range("N103").Value = -6304.3334

# Row 113
$ws.Range("H113").Value = 3684.5908
$ws.Range("J113").Value = 4063.0667
$ws.Range("L113").Value = 4063.0667
$ws.Range("N113").Value = -10571.0667

# Row 132
$ws.Range("H132").Value = 5147.604
$ws.Range("I132").Value = 3860.081
$ws.Range("K132").Value = 11580.243
$ws.Range("M132").Value = -9050.243

# ----- ARM sheet -----
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 5816.4097
$ws.Range("I32").Value = 5461.7124
$ws.Range("K32").Value = 5461.7124
$ws.Range("M32").Value = -5174.7124

# Row 132
$ws.Range("H132").Value = 3734.3704
$ws.Range("I132").Value = 1487.4722
$ws.Range("J132").Value = 8228.166999999999
$ws.Range("K132").Value = 4462.4166
$ws.Range("L132").Value = 24684.501
$ws.Range("M132").Value = -1932.4166
$ws.Range("N132").Value = -29744.501

# ----- BSM sheet -----
$ws = $wb.Worksheets.Item("BSM")

# Row 99
$ws.Range("H99").Value = 2174.4546
$ws.Range("I99").Value = 1679.8889
$ws.Range("K99").Value = 1679.8889
$ws.Range("M99").Value = -181.8888999999999

# Row 105
$ws.Range("H105").Value = 2520.0715
$ws.Range("I105").Value = 2360.9092
$ws.Range("K105").Value = 2360.9092
$ws.Range("M105").Value = -613.9092000000001

# Row 134
$ws.Range("H134").Value = 5276.1875
$ws.Range("I134").Value = 2607.2
$ws.Range("J134").Value = 7182.607
$ws.Range("K134").Value = 7821.599999999999
$ws.Range("L134").Value = 21547.821
$ws.Range("M134").Value = -5286.599999999999
$ws.Range("N134").Value = -26617.821

# ----- CRP sheet -----
$ws = $wb.Worksheets.Item("CRP")

# Row 132
$ws.Range("H132").Value = 2903.2
$ws.Range("I132").Value = 1910.9412
$ws.Range("J132").Value = 3840.3333
$ws.Range("K132").Value = 5732.8236
$ws.Range("L132").Value = 11520.9999
$ws.Range("M132").Value = -3202.8236
$ws.Range("N132").Value = -16580.9999

# ----- CUL sheet -----
$ws = $wb.Worksheets.Item("CUL")

# Row 3
$ws.Range("H3").Value = 6109.091
$ws.Range("I3").Value = 4733.3335
$ws.Range("J3").Value = 7760
$ws.Range("K3").Value = 14200.0005
$ws.Range("L3").Value = 23280
$ws.Range("M3").Value = -14088.0005
$ws.Range("N3").Value = -23504

# Row 18
$ws.Range("H18").Value = 1350.7222
$ws.Range("I18").Value = 406.66666
$ws.Range("J18").Value = 3238.8333
$ws.Range("K18").Value = 1219.99998
$ws.Range("L18").Value = 9716.499899999999
$ws.Range("M18").Value = -1050.99998
$ws.Range("N18").Value = -10054.4999

# Row 107
$ws.Range("H107").Value = 27778546
$ws.Range("I107").Value = 250000140
$ws.Range("J107").Value = 847.75
$ws.Range("K107").Value = 750000420
$ws.Range("L107").Value = 2543.25
$ws.Range("M107").Value = -749998500
$ws.Range("N107").Value = -6383.25

# Row 113
$ws.Range("H113").Value = 703.17645
$ws.Range("J113").Value = 841.7646999999999
$ws.Range("L113").Value = 2525.2941
$ws.Range("N113").Value = -6865.2941

# Row 115
$ws.Range("H115").Value = 2898.3333
$ws.Range("J115").Value = 3655.1667
$ws.Range("L115").Value = 10965.5001
$ws.Range("N115").Value = -13315.5001

# Row 133
$ws.Range("H133").Value = 3050.818
$ws.Range("I133").Value = 1351.8
$ws.Range("J133").Value = 4466.6665
$ws.Range("K133").Value = 4055.4
$ws.Range("L133").Value = 13399.9995
$ws.Range("M133").Value = 1004.6
$ws.Range("N133").Value = -23519.9995

# Row 134
$ws.Range("H134").Value = 3048.3684
$ws.Range("I134").Value = 1591.9
$ws.Range("J134").Value = 4666.6665
$ws.Range("K134").Value = 4775.700000000001
$ws.Range("L134").Value = 13999.9995
$ws.Range("M134").Value = 294.2999999999993
$ws.Range("N134").Value = -24139.9995

# ----- GSM sheet -----
$ws = $wb.Worksheets.Item("GSM")

# Row 113
$ws.Range("H113").Value = 3283.3333
$ws.Range("I113").Value = 3283.3333
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3283.3333
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1113.3333
$ws.Range("N113").ClearContents()

# Row 124
$ws.Range("H124").Value = 48000
$ws.Range("J124").Value = 48000
$ws.Range("L124").Value = 48000
$ws.Range("N124").Value = -57820

# Row 131
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# Row 132
$ws.Range("H132").Value = 2858.862
$ws.Range("I132").Value = 2426.9
$ws.Range("J132").Value = 3086.2104
$ws.Range("K132").Value = 7280.700000000001
$ws.Range("L132").Value = 9258.6312
$ws.Range("M132").Value = -4750.700000000001
$ws.Range("N132").Value = -14318.6312

# Row 134
$ws.Range("H134").Value = 28665.2
$ws.Range("J134").Value = 28665.2
$ws.Range("L134").Value = 85995.60000000001
$ws.Range("N134").Value = -91065.60000000001

# Row 136
$ws.Range("H136").Value = 23502.777
$ws.Range("J136").Value = 23502.777
$ws.Range("L136").Value = 70508.33099999999
$ws.Range("N136").Value = -75608.33099999999

# ----- LTW sheet -----
$ws = $wb.Worksheets.Item("LTW")

# Row 16
$ws.Range("H16").Value = 2355.7144
$ws.Range("I16").Value = 1493.6364
$ws.Range("J16").Value = 5516.6665
$ws.Range("K16").Value = 1493.6364
$ws.Range("L16").Value = 5516.6665
$ws.Range("M16").Value = -1323.6364
$ws.Range("N16").Value = -5856.6665

# Row 40
$ws.Range("H40").Value = 142862900
$ws.Range("I40").Value = 250002820
$ws.Range("J40").Value = 9666.666999999999
$ws.Range("K40").Value = 250002820
$ws.Range("L40").Value = 9666.666999999999
$ws.Range("M40").Value = -250002684
$ws.Range("N40").Value = -9938.666999999999

# Row 64
$ws.Range("H64").Value = 19000
$ws.Range("J64").Value = 19000
$ws.Range("L64").Value = 19000
$ws.Range("N64").Value = -19450

# Row 67
$ws.Range("H67").Value = 19000
$ws.Range("J67").Value = 19000
$ws.Range("L67").Value = 19000
$ws.Range("N67").Value = -20560

# ----- WVR sheet -----
$ws = $wb.Worksheets.Item("WVR")

# Row 63
$ws.Range("H63").Value = 20000
$ws.Range("J63").Value = 20000
$ws.Range("L63").Value = 20000
$ws.Range("N63").Value = -21248

# Row 66
$ws.Range("H66").Value = 20000
$ws.Range("J66").Value = 20000
$ws.Range("L66").Value = 60000
$ws.Range("N66").Value = -66240

# Row 96
$ws.Range("H96").Value = 6332.2144
$ws.Range("I96").Value = 2560.8333
$ws.Range("J96").Value = 9160.75
$ws.Range("K96").Value = 2560.8333
$ws.Range("L96").Value = 9160.75
$ws.Range("M96").Value = -1187.8333
$ws.Range("N96").Value = -11906.75

# Row 107
$ws.Range("H107").Value = 833.3333
$ws.Range("I107").Value = 725
$ws.Range("J107").Value = 920
$ws.Range("K107").Value = 2175
$ws.Range("L107").Value = 2760
$ws.Range("M107").Value = -255
$ws.Range("N107").Value = -6600
